$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# 1) First three rows' values become "0M"
$t.Cell(1, 1).Range.Text = "0M"
$t.Cell(2, 1).Range.Text = "0M"
$t.Cell(3, 1).Range.Text = "0M"

# 2) Insert 10 new single-value rows right after row 3 (i.e. before current row 4).
#    Rows.Add(refRow) always inserts immediately before refRow, so to end up with
#    the values in the intended ascending order we add them back-to-front.
$newValues = @("496", "0.00002", "0.00010", "0.00004", "0.00001", "0.00004", "0.00005", "0.00005", "0.01976", "100.0")
$refRow = $t.Rows.Item(4)
for ($i = $newValues.Length - 1; $i -ge 0; $i--) {
    $newRow = $t.Rows.Add($refRow)
    $newRow.Cells.Item(1).Range.Text = $newValues[$i]
    $refRow = $newRow
}

# 3) Collapse the three tab-separated multi-run rows (now shifted down by 10 rows,
#    i.e. originally rows 34/35/36, now rows 44/45/46) into single plain values.
$t.Cell(44, 1).Range.Text = "99.99"
$t.Cell(45, 1).Range.Text = "0.02"
$t.Cell(46, 1).Range.Text = "150"
